# Update the cryptos price list: refreshed Price figures (column D), plus a
# re-ranking of rows 18-24 (coin drops to top, others shift down one slot)
# and a couple of "Bestin24h"/label tweaks on the Volume(1h) column (E).
#
# Price cells hold their numbers as literal text (e.g. "245.48", keeping
# leading/trailing zeros like "0.03070"), so each one is stamped with a
# Text number format before the write - otherwise Excel would silently
# coerce the string into a numeric value and drop formatting like that.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.48"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.03"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05967"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.387"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8123"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9674"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1428"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07414"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03494"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03070"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09401"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.996"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001603"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04796"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005915"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006261"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004139"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009857"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00009704"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.741"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.166"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03906"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.003038"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1074"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002701"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.005366"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005317"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.03922"
